$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68; this shifts existing rows 68-144 down to 69-145,
# preserving their values and formatting (including the date style on column D).
$ws.Rows(68).Insert()

# Populate the newly inserted row 68 with the new record's data.
$ws.Cells.Item(68, 1).Value = 10
$ws.Cells.Item(68, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(68, 3).Value = "La Araucanía"
$ws.Cells.Item(68, 4).Value = 44482
$ws.Cells.Item(68, 5).Value = 9
$ws.Cells.Item(68, 6).Value = 100112052
$ws.Cells.Item(68, 7).Value = "Albahaca"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 20
$ws.Cells.Item(68, 11).Value = 7000
$ws.Cells.Item(68, 12).Value = 7000
$ws.Cells.Item(68, 13).Value = 7000
$ws.Cells.Item(68, 14).Value = "$/paquete"
$ws.Cells.Item(68, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(68, 16).Value = 7000
$ws.Cells.Item(68, 17).Value = 1
$ws.Cells.Item(68, 18).Value = "Hortaliza"
